$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host ("Sheet name: " + $ws.Name)
$lo = $ws.ListObjects.Item(1)
Write-Host ("Table name: " + $lo.Name)
Write-Host ("Table range: " + $lo.Range.Address())
Write-Host ("ListRows count: " + $lo.ListRows.Count)
